# Applies the my_Events.xlsx update:
#  - Inserts a new "ID" column (A) ahead of Day/Month/Year/Description,
#    shifting the existing columns right by one.
#  - Gives the new ID header (A1) an underlined font.
#  - Clears the old "Birthday Tom" row, leaving row 2 blank except for a
#    formatted (scientific-number) A2 cell.
#  - Rewrites the event rows (3-8) with unique numeric IDs plus the
#    day/month/year/description values, and appends three new events.
#  - Adjusts column widths and the active selection / page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Apply the two "real" cell styles first so they mint style indexes
#    1 and 2 (matching a from-scratch authoring order). The later
#    text-coercion trick (step 3) mints one more throwaway style index
#    that no cell ends up referencing.
# ---------------------------------------------------------------------
$ws.Range("A1").Font.Underline = $true
$ws.Range("A2").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------
# 2) Header row + shifted columns.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Day"
$ws.Range("C1").Value = "Month"
$ws.Range("D1").Value = "Year"
$ws.Range("E1").Value = "Description"

# ---------------------------------------------------------------------
# 3) Data rows. Day/Month/Year and the ID values all look like plain
#    numbers, so Excel would silently coerce them (dropping leading
#    zeros, and losing precision on the 17-20 digit IDs) unless the
#    cell is first marked as text. Prefixing with an apostrophe forces
#    text entry; ClearFormats() afterwards drops the resulting
#    "quotePrefix" flag again so the cells end up back at the default
#    (unstyled) cell format while keeping their text content/type.
# ---------------------------------------------------------------------
$data = @(
  @{ Row = 3; Id = "15293965254512808429"; Day = "05"; Month = "12"; Year = "2022"; Desc = "Test2" },
  @{ Row = 4; Id = "17779186155504800237"; Day = "01"; Month = "02"; Year = "2022"; Desc = "Test1" },
  @{ Row = 5; Id = "18073954687356768749"; Day = "04"; Month = "01"; Year = "2022"; Desc = "Test2" },
  @{ Row = 6; Id = "13550180417212909";    Day = "12"; Month = "12"; Year = "2023"; Desc = "Test3" },
  @{ Row = 7; Id = "504221220135834093";   Day = "29"; Month = "04"; Year = "2022"; Desc = "Test4" },
  @{ Row = 8; Id = "1692833795598258669";  Day = "30"; Month = "10"; Year = "2024"; Desc = "Test5" }
)

$textRanges = @()
foreach ($d in $data) {
  $r = $d.Row
  $ws.Range("A$r").Value = "'" + $d.Id
  $ws.Range("B$r").Value = "'" + $d.Day
  $ws.Range("C$r").Value = "'" + $d.Month
  $ws.Range("D$r").Value = "'" + $d.Year
  $ws.Range("E$r").Value = $d.Desc
  $textRanges += "A$r"
  $textRanges += "B$r"
  $textRanges += "C$r"
  $textRanges += "D$r"
}

foreach ($addr in $textRanges) {
  $ws.Range($addr).ClearFormats()
}

# ---------------------------------------------------------------------
# 4) Row 2: drop the old "Birthday Tom" event entirely, leaving only
#    the formatted-but-empty A2 cell behind.
# ---------------------------------------------------------------------
$ws.Range("B2:E2").ClearContents()
$ws.Range("A2").ClearContents()

# ---------------------------------------------------------------------
# 5) Column widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.7
$ws.Columns.Item(4).ColumnWidth = 8.45
$ws.Columns.Item(5).ColumnWidth = 10.7

# ---------------------------------------------------------------------
# 6) Selection + page setup.
# ---------------------------------------------------------------------
$ws.Range("L5").Select() | Out-Null
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9
